# Updates the cryptocurrency price/volume table on the active worksheet
# to reflect freshly scraped values (GitHub Actions data refresh).
#
# Each entry below describes a single cell edit: the cell reference, the
# new text value, and whether the value must be forced to stay as TEXT
# (because it looks numeric, e.g. "1.00" or "0.994", and Excel would
# otherwise silently coerce it into a number and drop formatting such as
# trailing/leading zeros).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '65.783.86'; ForceText = 0 },
    @{ Cell = 'E2'; Value = '  +4.98%  '; ForceText = 0 },
    @{ Cell = 'D3'; Value = '2.622.63'; ForceText = 0 },
    @{ Cell = 'E3'; Value = '  +6.80%  '; ForceText = 0 },
    @{ Cell = 'E4'; Value = '  -0.01%  '; ForceText = 0 },
    @{ Cell = 'D5'; Value = '589.21'; ForceText = 1 },
    @{ Cell = 'E5'; Value = '  +3.22%  '; ForceText = 0 },
    @{ Cell = 'D6'; Value = '155.23'; ForceText = 1 },
    @{ Cell = 'E6'; Value = '  +6.27%  '; ForceText = 0 },
    @{ Cell = 'D7'; Value = '1.00'; ForceText = 1 },
    @{ Cell = 'E7'; Value = '  +0.04%  '; ForceText = 0 },
    @{ Cell = 'E8'; Value = '  +3.15%  '; ForceText = 0 },
    @{ Cell = 'D9'; Value = '2.621.72'; ForceText = 0 },
    @{ Cell = 'E9'; Value = '  +6.75%  '; ForceText = 0 },
    @{ Cell = 'E10'; Value = '  +4.06%  '; ForceText = 0 },
    @{ Cell = 'E11'; Value = '  -1.80%  '; ForceText = 0 },
    @{ Cell = 'E12'; Value = '  +4.19%  '; ForceText = 0 },
    @{ Cell = 'D13'; Value = '5.32'; ForceText = 1 },
    @{ Cell = 'E13'; Value = '  +2.08%  '; ForceText = 0 },
    @{ Cell = 'D14'; Value = '29.33'; ForceText = 1 },
    @{ Cell = 'E14'; Value = '  +2.58%  '; ForceText = 0 },
    @{ Cell = 'D15'; Value = '3.098.56'; ForceText = 0 },
    @{ Cell = 'E15'; Value = '  +7.37%  '; ForceText = 0 },
    @{ Cell = 'D16'; Value = '0.0000182'; ForceText = 1 },
    @{ Cell = 'E16'; Value = '  +4.57%  '; ForceText = 0 },
    @{ Cell = 'D17'; Value = '65.654.57'; ForceText = 0 },
    @{ Cell = 'E17'; Value = '  +5.05%  '; ForceText = 0 },
    @{ Cell = 'D18'; Value = '2.613.38'; ForceText = 0 },
    @{ Cell = 'E18'; Value = '  +6.38%  '; ForceText = 0 },
    @{ Cell = 'D19'; Value = '8.26'; ForceText = 1 },
    @{ Cell = 'E19'; Value = '  +7.91%  '; ForceText = 0 },
    @{ Cell = 'D20'; Value = '11.22'; ForceText = 1 },
    @{ Cell = 'E20'; Value = '  +4.36%  '; ForceText = 0 },
    @{ Cell = 'D21'; Value = '355.80'; ForceText = 1 },
    @{ Cell = 'E21'; Value = '  +10.87%  '; ForceText = 0 },
    @{ Cell = 'E22'; Value = '  +4.28%  '; ForceText = 0 },
    @{ Cell = 'E23'; Value = '  +1.64%  '; ForceText = 0 },
    @{ Cell = 'D24'; Value = '1.00'; ForceText = 1 },
    @{ Cell = 'E24'; Value = '  -0.05%  '; ForceText = 0 },
    @{ Cell = 'D25'; Value = '9.97'; ForceText = 1 },
    @{ Cell = 'E25'; Value = '  +0.45%  '; ForceText = 0 },
    @{ Cell = 'D26'; Value = '66.28'; ForceText = 1 },
    @{ Cell = 'E26'; Value = '  +1.85%  '; ForceText = 0 },
    @{ Cell = 'D27'; Value = '632.71'; ForceText = 1 },
    @{ Cell = 'E27'; Value = '  -1.61%  '; ForceText = 0 },
    @{ Cell = 'D28'; Value = '0.0000106'; ForceText = 1 },
    @{ Cell = 'E28'; Value = '  +10.58%  '; ForceText = 0 },
    @{ Cell = 'D29'; Value = '2.729.43'; ForceText = 0 },
    @{ Cell = 'E29'; Value = '  +6.62%  '; ForceText = 0 },
    @{ Cell = 'E30'; Value = '  +5.93%  '; ForceText = 0 },
    @{ Cell = 'D31'; Value = '0.994'; ForceText = 1 },
    @{ Cell = 'E31'; Value = '  -0.19%  '; ForceText = 0 },
    @{ Cell = 'D32'; Value = '8.27'; ForceText = 1 },
    @{ Cell = 'E32'; Value = '  +5.92%  '; ForceText = 0 },
    @{ Cell = 'D33'; Value = '1.91'; ForceText = 1 },
    @{ Cell = 'E33'; Value = '  +5.48%  '; ForceText = 0 },
    @{ Cell = 'D34'; Value = '0.138'; ForceText = 1 },
    @{ Cell = 'E34'; Value = '  +4.43%  '; ForceText = 0 },
    @{ Cell = 'D35'; Value = '1.64'; ForceText = 1 },
    @{ Cell = 'E35'; Value = '  +9.58%  '; ForceText = 0 },
    @{ Cell = 'D36'; Value = '0.999'; ForceText = 1 },
    @{ Cell = 'E36'; Value = '  +0.09%  '; ForceText = 0 },
    @{ Cell = 'D37'; Value = '4.98'; ForceText = 1 },
    @{ Cell = 'E37'; Value = '  +7.26%  '; ForceText = 0 },
    @{ Cell = 'E38'; Value = '  +5.58%  '; ForceText = 0 },
    @{ Cell = 'D39'; Value = '19.41'; ForceText = 1 },
    @{ Cell = 'E39'; Value = '  +4.94%  '; ForceText = 0 },
    @{ Cell = 'D40'; Value = '2.90'; ForceText = 1 },
    @{ Cell = 'E40'; Value = '  +6.50%  '; ForceText = 0 },
    @{ Cell = 'D41'; Value = '155.77'; ForceText = 1 },
    @{ Cell = 'E41'; Value = '  +3.56%  '; ForceText = 0 },
    @{ Cell = 'E42'; Value = '  +2.92%  '; ForceText = 0 },
    @{ Cell = 'E43'; Value = '  +6.48%  '; ForceText = 0 },
    @{ Cell = 'B44'; Value = 'OKB'; ForceText = 0 },
    @{ Cell = 'C44'; Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'; ForceText = 0 },
    @{ Cell = 'D44'; Value = '42.18'; ForceText = 1 },
    @{ Cell = 'E44'; Value = '  +1.21%  '; ForceText = 0 },
    @{ Cell = 'B45'; Value = 'Aave'; ForceText = 0 },
    @{ Cell = 'C45'; Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'; ForceText = 0 },
    @{ Cell = 'D45'; Value = '163.37'; ForceText = 1 },
    @{ Cell = 'E45'; Value = '  +6.86%  '; ForceText = 0 },
    @{ Cell = 'B46'; Value = 'USDe'; ForceText = 0 },
    @{ Cell = 'C46'; Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'; ForceText = 0 },
    @{ Cell = 'D46'; Value = '0.999'; ForceText = 1 },
    @{ Cell = 'E46'; Value = '  -0.02%  '; ForceText = 0 },
    @{ Cell = 'B47'; Value = 'WhiteBITCoin'; ForceText = 0 },
    @{ Cell = 'C47'; Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'; ForceText = 0 },
    @{ Cell = 'D47'; Value = '16.20'; ForceText = 1 },
    @{ Cell = 'E47'; Value = '  +5.26%  '; ForceText = 0 },
    @{ Cell = 'B48'; Value = 'Filecoin'; ForceText = 0 },
    @{ Cell = 'C48'; Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'; ForceText = 0 },
    @{ Cell = 'D48'; Value = '3.77'; ForceText = 1 },
    @{ Cell = 'E48'; Value = '  +6.48%  '; ForceText = 0 },
    @{ Cell = 'B49'; Value = 'InjectiveProtocol'; ForceText = 0 },
    @{ Cell = 'C49'; Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'; ForceText = 0 },
    @{ Cell = 'D49'; Value = '21.78'; ForceText = 1 },
    @{ Cell = 'E49'; Value = '  +8.15%  '; ForceText = 0 },
    @{ Cell = 'B50'; Value = 'Mantle'; ForceText = 0 },
    @{ Cell = 'C50'; Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'; ForceText = 0 },
    @{ Cell = 'D50'; Value = '0.637'; ForceText = 1 },
    @{ Cell = 'E50'; Value = '  +5.39%  '; ForceText = 0 },
    @{ Cell = 'B51'; Value = 'Hedera'; ForceText = 0 },
    @{ Cell = 'C51'; Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'; ForceText = 0 },
    @{ Cell = 'D51'; Value = '0.0530'; ForceText = 1 },
    @{ Cell = 'E51'; Value = '  +5.19%  '; ForceText = 0 }
)

foreach ($update in $updates) {
    $range = $ws.Range($update.Cell)
    if ($update.ForceText -eq 1) {
        # Force a text number format so numeric-looking strings (e.g. "1.00",
        # "0.0000182") are written verbatim as text instead of being
        # normalized into numbers.
        $range.NumberFormat = "@"
        $range.Value = $update.Value
        # Restore the default "Normal" style so no stray formatting/style
        # index is left behind on the cell.
        $range.Style = "Normal"
    } else {
        $range.Value = $update.Value
    }
}
